$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 408: S33_G01_TB001
$ws.Rows.Item(408).RowHeight = 58.3
$ws.Range("A408").Value2 = 'S33'
$ws.Range("B408").Value2 = 'G01'
$ws.Range("C408").Value2 = 'Holdings with intent (Goal View) MVP'
$ws.Range("D408").Value2 = 'S33_G01_TB001'
$ws.Range("E408").Value2 = 'Backend (MIGRATION): Add holding_goals persistence (label, review_date, target_type, target_value, note, created_at, updated_at) keyed by user+broker+symbol+exchange.'
$ws.Range("G408").Value2 = 'implemented'
$ws.Range("H408").Value2 = 'Holding goals table + model added (label, review date, target, note, timestamps).'

# Row 409: S33_G01_TB002
$ws.Rows.Item(409).RowHeight = 58.3
$ws.Range("A409").Value2 = 'S33'
$ws.Range("B409").Value2 = 'G01'
$ws.Range("C409").Value2 = 'Holdings with intent (Goal View) MVP'
$ws.Range("D409").Value2 = 'S33_G01_TB002'
$ws.Range("E409").Value2 = 'Backend: Add goals CRUD (list by broker, upsert by symbol/exchange, delete) with label-based default review date and single target type.'
$ws.Range("G409").Value2 = 'implemented'
$ws.Range("H409").Value2 = 'Goals CRUD endpoints added with defaults + validation.'

# Row 410: S33_G01_TF001
$ws.Rows.Item(410).RowHeight = 58.3
$ws.Range("A410").Value2 = 'S33'
$ws.Range("B410").Value2 = 'G01'
$ws.Range("C410").Value2 = 'Holdings with intent (Goal View) MVP'
$ws.Range("D410").Value2 = 'S33_G01_TF001'
$ws.Range("E410").Value2 = 'Frontend: Add Goal View toggle in Holdings, new columns (label, review date, days, status, target, away %, note) and filters (All/Overdue/Due Soon/Near Target/Missing).'
$ws.Range("G410").Value2 = 'implemented'
$ws.Range("H410").Value2 = 'Goal View toggle + columns + filters wired in Holdings.'

# Row 411: S33_G01_TF002
$ws.Rows.Item(411).RowHeight = 58.3
$ws.Range("A411").Value2 = 'S33'
$ws.Range("B411").Value2 = 'G01'
$ws.Range("C411").Value2 = 'Holdings with intent (Goal View) MVP'
$ws.Range("D411").Value2 = 'S33_G01_TF002'
$ws.Range("E411").Value2 = 'Frontend: Implement Edit Goal drawer (fixed labels, visible defaults, single target type, note) with save/update to goals API.'
$ws.Range("G411").Value2 = 'implemented'
$ws.Range("H411").Value2 = 'Goal edit dialog implemented with label defaults + target preview.'

# Row 412: S33_G01_TF003
$ws.Rows.Item(412).RowHeight = 58.3
$ws.Range("A412").Value2 = 'S33'
$ws.Range("B412").Value2 = 'G01'
$ws.Range("C412").Value2 = 'Holdings with intent (Goal View) MVP'
$ws.Range("D412").Value2 = 'S33_G01_TF003'
$ws.Range("E412").Value2 = 'Frontend: Add soft prompts for missing goals (badge, CTA "Set missing goals", Missing filter; no hard block).'
$ws.Range("G412").Value2 = 'implemented'
$ws.Range("H412").Value2 = 'Soft prompts for missing goals + filter/CTA added.'

# Row 413: S33_G01_TD001
$ws.Rows.Item(413).RowHeight = 58.3
$ws.Range("A413").Value2 = 'S33'
$ws.Range("B413").Value2 = 'G01'
$ws.Range("C413").Value2 = 'Holdings with intent (Goal View) MVP'
$ws.Range("D413").Value2 = 'S33_G01_TD001'
$ws.Range("E413").Value2 = 'Docs/QA: Add manual QA checklist for Goal View MVP (missing goals, due soon, overdue, near target, edit save).'
$ws.Range("G413").Value2 = 'implemented'
$ws.Range("H413").Value2 = 'QA checklist added in docs/qa/holdings_goal_view.md'

# Row 414: S33_G02_TB001
$ws.Rows.Item(414).RowHeight = 58.3
$ws.Range("A414").Value2 = 'S33'
$ws.Range("B414").Value2 = 'G02'
$ws.Range("C414").Value2 = 'Goal View v1: CSV import + presets'
$ws.Range("D414").Value2 = 'S33_G02_TB001'
$ws.Range("E414").Value2 = 'Backend: Add bulk import endpoint for holding goals with symbol normalization (NSE:/BSE:), match summary, and per-row errors.'
$ws.Range("G414").Value2 = 'planned'
$ws.Range("H414").Value2 = 'Est: 3h | Area: BE | Paths: backend/app/api/holdings_goals.py, backend/app/services/holdings_goals.py | Deps: S33_G01_TB002'

# Row 415: S33_G02_TB002
$ws.Rows.Item(415).RowHeight = 58.3
$ws.Range("A415").Value2 = 'S33'
$ws.Range("B415").Value2 = 'G02'
$ws.Range("C415").Value2 = 'Goal View v1: CSV import + presets'
$ws.Range("D415").Value2 = 'S33_G02_TB002'
$ws.Range("E415").Value2 = 'Backend: Add CSV mapping preset CRUD (save, list, delete) scoped to user.'
$ws.Range("G415").Value2 = 'planned'
$ws.Range("H415").Value2 = 'Est: 2h | Area: BE | Paths: backend/app/api/holdings_goals.py, backend/app/models/holdings.py | Deps: S33_G01_TB001'

# Row 416: S33_G02_TF001
$ws.Rows.Item(416).RowHeight = 58.3
$ws.Range("A416").Value2 = 'S33'
$ws.Range("B416").Value2 = 'G02'
$ws.Range("C416").Value2 = 'Goal View v1: CSV import + presets'
$ws.Range("D416").Value2 = 'S33_G02_TF001'
$ws.Range("E416").Value2 = 'Frontend: Build CSV import wizard (upload, preview, column mapping, preset save, import summary) for Goal View.'
$ws.Range("G416").Value2 = 'planned'
$ws.Range("H416").Value2 = 'Est: 4h | Area: FE | Paths: frontend/src/views/HoldingsPage.tsx, frontend/src/components/GoalImportDialog.tsx (new), frontend/src/services/holdingsGoals.ts | Deps: S33_G02_TB001, S33_G02_TB002'

# Row 417: S33_G02_TD001
$ws.Rows.Item(417).RowHeight = 58.3
$ws.Range("A417").Value2 = 'S33'
$ws.Range("B417").Value2 = 'G02'
$ws.Range("C417").Value2 = 'Goal View v1: CSV import + presets'
$ws.Range("D417").Value2 = 'S33_G02_TD001'
$ws.Range("E417").Value2 = 'Docs/QA: CSV import checklist (mapping, presets, unmatched symbols, update counts).'
$ws.Range("G417").Value2 = 'planned'
$ws.Range("H417").Value2 = 'Est: 1h | Area: DOCS | Paths: docs/qa/holdings_goal_import.md (new) | Deps: S33_G02_TF001'

# Row 418: S33_G03_TB001
$ws.Rows.Item(418).RowHeight = 58.3
$ws.Range("A418").Value2 = 'S33'
$ws.Range("B418").Value2 = 'G03'
$ws.Range("C418").Value2 = 'Goal View v2: alerts + review workflow'
$ws.Range("D418").Value2 = 'S33_G03_TB001'
$ws.Range("E418").Value2 = 'Backend: Add review workflow support (last_reviewed_at/history) and stop target alerts after review date unless extended.'
$ws.Range("G418").Value2 = 'planned'
$ws.Range("H418").Value2 = 'Est: 3h | Area: BE | Paths: backend/app/models/holdings.py, backend/app/services/alerts_v3.py | Deps: S33_G01_TB001'

# Row 419: S33_G03_TF001
$ws.Rows.Item(419).RowHeight = 58.3
$ws.Range("A419").Value2 = 'S33'
$ws.Range("B419").Value2 = 'G03'
$ws.Range("C419").Value2 = 'Goal View v2: alerts + review workflow'
$ws.Range("D419").Value2 = 'S33_G03_TF001'
$ws.Range("E419").Value2 = 'Frontend: Add review actions (extend, snooze) and review history panel in Goal View.'
$ws.Range("G419").Value2 = 'planned'
$ws.Range("H419").Value2 = 'Est: 3h | Area: FE | Paths: frontend/src/views/HoldingsPage.tsx, frontend/src/components/GoalReviewPanel.tsx (new) | Deps: S33_G03_TB001'

# Row 420: S33_G03_TD001
$ws.Rows.Item(420).RowHeight = 58.3
$ws.Range("A420").Value2 = 'S33'
$ws.Range("B420").Value2 = 'G03'
$ws.Range("C420").Value2 = 'Goal View v2: alerts + review workflow'
$ws.Range("D420").Value2 = 'S33_G03_TD001'
$ws.Range("E420").Value2 = 'Docs: Goal View rollout notes + guardrails (no auto-sell, intent-first behavior).'
$ws.Range("G420").Value2 = 'planned'
$ws.Range("H420").Value2 = 'Est: 1h | Area: DOCS | Paths: docs/holdings_goal_view.md (new) | Deps: S33_G03_TF001'

